# Generate Report for handback
# File "7593a10e-691d-4a59-b935-bd2d1ef4e50c.md" has been handed back for both
# locales (zh-cn, de-de) and is now "in sync with en-US". Update the Overview
# sheet status, the per-locale sheets' Status + Latest Handback DateTime.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: row for 7593a10e-... (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet: row for 7593a10e-... (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $statusHandedBack
$wsZhCn.Range("G3").Value = "2016-01-13 15:49:27"

# --- de-de sheet: row for 7593a10e-... (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $statusHandedBack
$wsDeDe.Range("G3").Value = "2016-01-13 15:49:46"
